# Apply updated crypto price/volume data per commit "Updated cryptos list on Sat Jun 17 09:55:58 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.594.03"
$ws.Range("E2").Value = "  +3.99%  "
# Row 3
$ws.Range("D3").Value = "1.743.26"
$ws.Range("E3").Value = "  +4.46%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").Formula = "'246.15"
# Row 6
$ws.Range("E6").Value = "  -0.01%  "
# Row 7
$ws.Range("D7").Formula = "'0.4822"
$ws.Range("E7").Value = "  +0.95%  "
# Row 8
$ws.Range("D8").Formula = "'0.2701"
$ws.Range("E8").Value = "  +3.48%  "
# Row 9
$ws.Range("D9").Formula = "'0.06264"
$ws.Range("E9").Value = "  +1.52%  "
# Row 10
$ws.Range("D10").Value = "1.744.40"
$ws.Range("E10").Value = "  +4.46%  "
# Row 11
$ws.Range("D11").Formula = "'0.07128"
$ws.Range("E11").Value = "  +1.91%  "
# Row 12
$ws.Range("E12").Value = "  +7.02%  "
# Row 13
$ws.Range("D13").Formula = "'0.6255"
$ws.Range("E13").Value = "  +6.04%  "
# Row 14
$ws.Range("D14").Formula = "'4.525"
$ws.Range("E14").Value = "  +3.39%  "
# Row 15
$ws.Range("D15").Formula = "'77.55"
$ws.Range("E15").Value = "  +2.88%  "
# Row 16
$ws.Range("D16").Formula = "'1.0000"
$ws.Range("E16").Value = "  -0.01%  "
# Row 17
$ws.Range("D17").Value = "26.594.67"
$ws.Range("E17").Value = "  +4.01%  "
# Row 18
$ws.Range("D18").Formula = "'1.000"
# Row 19
$ws.Range("D19").Formula = "'0.000006904"
$ws.Range("E19").Value = "  +2.47%  "
# Row 20
$ws.Range("D20").Formula = "'11.74"
$ws.Range("E20").Value = "  +2.80%  "
# Row 21
$ws.Range("D21").Value = "1.968.32"
$ws.Range("E21").Value = "  +4.40%  "
# Row 22
$ws.Range("D22").Formula = "'4.632"
$ws.Range("E22").Value = "  +4.22%  "
# Row 23
$ws.Range("D23").Formula = "'8.850"
$ws.Range("E23").Value = "  +0.60%  "
# Row 24
$ws.Range("D24").Formula = "'5.376"
$ws.Range("E24").Value = "  +2.17%  "
# Row 25
$ws.Range("D25").Formula = "'136.23"
$ws.Range("E25").Value = "  -0.49%  "
# Row 26
$ws.Range("D26").Formula = "'15.40"
$ws.Range("E26").Value = "  +2.56%  "
# Row 27
$ws.Range("E27").Value = "  +5.91%  "
# Row 28
$ws.Range("D28").Formula = "'1.431"
$ws.Range("E28").Value = "  +3.40%  "
# Row 29
$ws.Range("D29").Formula = "'107.02"
$ws.Range("E29").Value = "  +2.19%  "
# Row 30
$ws.Range("D30").Formula = "'4.008"
$ws.Range("E30").Value = "  +0.24%  "
# Row 31
$ws.Range("D31").Formula = "'3.751"
$ws.Range("E31").Value = "  +3.32%  "
# Row 32
$ws.Range("D32").Formula = "'0.07891"
$ws.Range("E32").Value = "  +0.37%  "
# Row 33
$ws.Range("D33").Formula = "'0.04622"
$ws.Range("E33").Value = "  +7.12%  "
# Row 34
$ws.Range("D34").Formula = "'2.618"
$ws.Range("E34").Value = "  -0.20%  "
# Row 35
$ws.Range("D35").Formula = "'0.6419"
$ws.Range("E35").Value = "  +6.28%  "
# Row 36
$ws.Range("D36").Formula = "'0.9991"
$ws.Range("E36").Value = "  +4.53%  "
# Row 37
$ws.Range("D37").Formula = "'0.9421"
$ws.Range("E37").Value = "  +0.52%  "
# Row 38
$ws.Range("D38").Formula = "'113.33"
$ws.Range("E38").Value = "  +15.64%  "
# Row 39
$ws.Range("E39").Value = "  +8.36%  "
# Row 40
$ws.Range("D40").Formula = "'2.427"
$ws.Range("E40").Value = "  -5.91%  "
# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Formula = "'1.002"
$ws.Range("E41").Value = "  +0.28%  "
# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Formula = "'5.775"
$ws.Range("E42").Value = "  +18.13%  "
# Row 43
$ws.Range("D43").Formula = "'0.01512"
# Row 44
$ws.Range("D44").Formula = "'0.3922"
$ws.Range("E44").Value = "  +4.30%  "
# Row 45
$ws.Range("D45").Formula = "'0.1224"
$ws.Range("E45").Value = "  +9.25%  "
# Row 46
$ws.Range("D46").Formula = "'6.737"
$ws.Range("E46").Value = "  +8.36%  "
# Row 47
$ws.Range("D47").Formula = "'0.05335"
$ws.Range("E47").Value = "  +1.31%  "
# Row 48
$ws.Range("D48").Formula = "'7.922"
$ws.Range("E48").Value = "  +6.17%  "
# Row 49
$ws.Range("D49").Formula = "'30.78"
$ws.Range("E49").Value = "  +2.68%  "
# Row 50
$ws.Range("D50").Formula = "'1.264"
$ws.Range("E50").Value = "  +4.83%  "
# Row 51
$ws.Range("D51").Formula = "'0.3456"
$ws.Range("E51").Value = "  +3.49%  "
